$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("A_SERIES_HOCHSCHUL",   "K_QUALI",  "Hochschulabschlüsse", "XXX"),
    @("A_SERIES_HOEHERQUAL",  "K_QUALI",  "Höherqualifizierende Berufsausbildung", "XXX"),
    @("A_SERIES_ZWEITAUSB",   "K_QUALI",  "Zweitausbildung, Ausbildung nach Abitur", "XXX"),
    @("A_SERIES_ERSTAUSB",    "K_QUALI",  "Erstausbildung", "XXX"),
    @("A_SERIES_FACHKR",      "K_SERIES", "Fachkräftepotenzial (insgesamt)", "XXX"),
    @("A_SERIES_QUALIFZIERT", "K_SERIES", "Akademisch Qualifizierte und beruflich Höherqualifizierte", "XXX")
)

$templateRow = 529
$startRow = 530

$ws.Range("A$templateRow`:D$templateRow").Copy()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row`:D$row").PasteSpecial(-4122)
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
